$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: new date (style like A11 / m-d-yyyy), activity, hours
$ws.Range("A11").Copy($ws.Range("A15"))
$ws.Range("A15").Value = 43152
$ws.Range("B15").Value = "Generic Dao/Testing"
$ws.Range("E15").Value = 2

# Row 16: new date (style like A12 / d-mmm), activity, hours
$ws.Range("A12").Copy($ws.Range("A16"))
$ws.Range("A16").Value = 43153
$ws.Range("B16").Value = "Fixing tests/DB"
$ws.Range("E16").Value = 2

# Row 17: new date (style like A11 / m-d-yyyy), activity, hours
$ws.Range("A11").Copy($ws.Range("A17"))
$ws.Range("A17").Value = 43157
$ws.Range("B17").Value = "AWS Setup / First Deploy"
$ws.Range("E17").Value = 6

# Row 18: new date only (style like A11 / m-d-yyyy)
$ws.Range("A11").Copy($ws.Range("A18"))
$ws.Range("A18").Value = 43158

# Selection moved to B18 in the saved view
$ws.Range("B18").Select() | Out-Null
